# Refresh the "ランサーズ" (Lancers) scraped-listing sheet with a newer
# snapshot: the capture timestamp moves from 01:27:18 to 06:23:43, most
# rows are replaced with the newer scrape's rows, and the trailing rows
# (8-18 in the old data) fall out of the refreshed window entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Drop every existing hyperlink (and their relationships) up front -
# rows are about to be rewritten/removed and new hyperlinks will be
# (re)created only for the rows that survive.
$ws.Hyperlinks.Delete()

# The old rows 8-18 no longer exist in the refreshed snapshot.
$ws.Rows("8:18").Delete()

# Row 2 keeps its content; only the capture timestamp advances.
$ws.Range("A2").Value = "2025-09-13 06:23:43"

# Row 3
$ws.Range("A3").Value = "2025-09-13 06:23:43"
$ws.Range("B3").Value = "【急募】アプリケーションAI駆動開発の仕上げをお手伝いください!"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5392937"
$ws.Range("G3").Value = 378
$ws.Range("H3").Value = "🔥AI,Ai ◆開発 ◇アプリ"

# Row 4
$ws.Range("A4").Value = "2025-09-13 06:23:43"
$ws.Range("B4").Value = "【急募】EA自動化システム構築の専門家を探しています!"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5392078"
$ws.Range("G4").Value = 98
$ws.Range("H4").Value = "◆自動化"

# Row 5
$ws.Range("A5").Value = "2025-09-13 06:23:43"
$ws.Range("B5").Value = "初回 WEB講習システムの開発"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5392840"
$ws.Range("G5").Value = 90
$ws.Range("H5").Value = "◆開発"

# Row 6 - note the new snapshot has no skill-summary (H) text for this row.
$ws.Range("A6").Value = "2025-09-13 06:23:43"
$ws.Range("B6").Value = "【集客支援】X(旧Twitter)、スレッズでの自動集客を実現したい!"
$ws.Range("C6").Value = "システム開発"
$ws.Range("D6").Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Range("E6").Value = "期限情報なし"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5392785"
$ws.Range("G6").Value = 10
$ws.Range("H6").ClearContents()

# Row 7 - also no skill-summary (H) text in the refreshed snapshot.
$ws.Range("A7").Value = "2025-09-13 06:23:43"
$ws.Range("B7").Value = "【AWSのプロ募集】事業成長を共に牽引するクラウドインフラの設計・構築パートナー募集中!"
$ws.Range("C7").Value = "システム開発"
$ws.Range("D7").Value = "~ 5,000 円 / 固定"
$ws.Range("E7").Value = "期限情報なし"
$ws.Range("F7").Value = "https://www.lancers.jp/work/detail/5392608"
$ws.Range("G7").Value = 10
$ws.Range("H7").ClearContents()

# Re-create the hyperlinks for the URL column on the rows that remain.
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5392661")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5392937")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5392078")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5392840")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5392785")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5392608")

# Column width tweaks (B: 47 -> 46, D: 32 -> 28). ColumnWidth is in
# "characters" and differs from the raw stored OOXML width by the fixed
# ~0.8333 char padding Excel adds, so back that out to land on the exact
# stored width the diff expects.
$ws.Columns("B").ColumnWidth = 46 - 0.8333333333333
$ws.Columns("D").ColumnWidth = 28 - 0.8333333333333
